# Weapon_Attack_Factors data rebuild: add base_round_time / minimum_round_time
# columns, the blunt_weapons / polearm_weapons weapon groups, and rename
# "hand_axe" -> "handaxe" (combat rework per commit message).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1
$ws.Range("A1").Value = "classification"
$ws.Range("B1").Value = "sub_category"
$ws.Range("C1").Value = "None"
$ws.Range("D1").Value = "cloth"
$ws.Range("E1").Value = "leather"
$ws.Range("F1").Value = "scale"
$ws.Range("G1").Value = "chain"
$ws.Range("H1").Value = "plate"
$ws.Range("I1").Value = "base_round_time"
$ws.Range("J1").Value = "minimum_round_time"

# Row 2
$ws.Range("A2").Value = "None"
$ws.Range("B2").Value = "None"
$ws.Range("C2").Value = 10
$ws.Range("D2").Value = 5
$ws.Range("E2").Value = 0
$ws.Range("F2").Value = -10
$ws.Range("G2").Value = -15
$ws.Range("H2").Value = -20
$ws.Range("I2").Value = 1
$ws.Range("J2").Value = 3

# Row 3
$ws.Range("A3").Value = "dagger"
$ws.Range("B3").Value = "edged_weapons"
$ws.Range("C3").Value = 27
$ws.Range("D3").Value = 25
$ws.Range("E3").Value = 20
$ws.Range("F3").Value = 9
$ws.Range("G3").Value = -2
$ws.Range("H3").Value = -18
$ws.Range("I3").Value = 1
$ws.Range("J3").Value = 3

# Row 4
$ws.Range("A4").Value = "rapier"
$ws.Range("B4").Value = "edged_weapons"
$ws.Range("C4").Value = 45
$ws.Range("D4").Value = 45
$ws.Range("E4").Value = 37
$ws.Range("F4").Value = 24
$ws.Range("G4").Value = 23
$ws.Range("H4").Value = -3
$ws.Range("I4").Value = 2
$ws.Range("J4").Value = 4

# Row 5
$ws.Range("A5").Value = "short_sword"
$ws.Range("B5").Value = "edged_weapons"
$ws.Range("C5").Value = 40
$ws.Range("D5").Value = 40
$ws.Range("E5").Value = 33
$ws.Range("F5").Value = 24
$ws.Range("G5").Value = 13
$ws.Range("H5").Value = 7
$ws.Range("I5").Value = 3
$ws.Range("J5").Value = 4

# Row 6
$ws.Range("A6").Value = "longsword"
$ws.Range("B6").Value = "edged_weapons"
$ws.Range("C6").Value = 41
$ws.Range("D6").Value = 41
$ws.Range("E6").Value = 39
$ws.Range("F6").Value = 37
$ws.Range("G6").Value = 25
$ws.Range("H6").Value = 17
$ws.Range("I6").Value = 4
$ws.Range("J6").Value = 4

# Row 7
$ws.Range("A7").Value = "broadsword"
$ws.Range("B7").Value = "edged_weapons"
$ws.Range("C7").Value = 36
$ws.Range("D7").Value = 36
$ws.Range("E7").Value = 33
$ws.Range("F7").Value = 30
$ws.Range("G7").Value = 25
$ws.Range("H7").Value = 18
$ws.Range("I7").Value = 5
$ws.Range("J7").Value = 5

# Row 8
$ws.Range("A8").Value = "claymore"
$ws.Range("B8").Value = "edged_weapons"
$ws.Range("C8").Value = 30
$ws.Range("D8").Value = 30
$ws.Range("E8").Value = 28
$ws.Range("F8").Value = 25
$ws.Range("G8").Value = 20
$ws.Range("H8").Value = 13
$ws.Range("I8").Value = 8
$ws.Range("J8").Value = 5

# Row 9
$ws.Range("A9").Value = "handaxe"
$ws.Range("B9").Value = "edged_weapons"
$ws.Range("C9").Value = 30
$ws.Range("D9").Value = 30
$ws.Range("E9").Value = 29
$ws.Range("F9").Value = 32
$ws.Range("G9").Value = 29
$ws.Range("H9").Value = 23
$ws.Range("I9").Value = 4
$ws.Range("J9").Value = 4

# Row 10
$ws.Range("A10").Value = "battle_axe"
$ws.Range("B10").Value = "edged_weapons"
$ws.Range("C10").Value = 35
$ws.Range("D10").Value = 35
$ws.Range("E10").Value = 32
$ws.Range("F10").Value = 35
$ws.Range("G10").Value = 32
$ws.Range("H10").Value = 25
$ws.Range("I10").Value = 6
$ws.Range("J10").Value = 5

# Row 11
$ws.Range("A11").Value = "whip"
$ws.Range("B11").Value = "blunt_weapons"
$ws.Range("C11").Value = 35
$ws.Range("D11").Value = 35
$ws.Range("E11").Value = 23
$ws.Range("F11").Value = 16
$ws.Range("G11").Value = 17
$ws.Range("H11").Value = 3
$ws.Range("I11").Value = 2
$ws.Range("J11").Value = 3

# Row 12
$ws.Range("A12").Value = "cudgel"
$ws.Range("B12").Value = "blunt_weapons"
$ws.Range("C12").Value = 20
$ws.Range("D12").Value = 20
$ws.Range("E12").Value = 18
$ws.Range("F12").Value = 20
$ws.Range("G12").Value = 24
$ws.Range("H12").Value = 17
$ws.Range("I12").Value = 4
$ws.Range("J12").Value = 3

# Row 13
$ws.Range("A13").Value = "mace"
$ws.Range("B13").Value = "blunt_weapons"
$ws.Range("C13").Value = 31
$ws.Range("D13").Value = 31
$ws.Range("E13").Value = 30
$ws.Range("F13").Value = 31
$ws.Range("G13").Value = 34
$ws.Range("H13").Value = 24
$ws.Range("I13").Value = 4
$ws.Range("J13").Value = 4

# Row 14
$ws.Range("A14").Value = "morning_star"
$ws.Range("B14").Value = "blunt_weapons"
$ws.Range("C14").Value = 33
$ws.Range("D14").Value = 33
$ws.Range("E14").Value = 33
$ws.Range("F14").Value = 30
$ws.Range("G14").Value = 34
$ws.Range("H14").Value = 25
$ws.Range("I14").Value = 5
$ws.Range("J14").Value = 4

# Row 15
$ws.Range("A15").Value = "flail"
$ws.Range("B15").Value = "blunt_weapons"
$ws.Range("C15").Value = 15
$ws.Range("D15").Value = 15
$ws.Range("E15").Value = 18
$ws.Range("F15").Value = 23
$ws.Range("G15").Value = 27
$ws.Range("H15").Value = 18
$ws.Range("I15").Value = 6
$ws.Range("J15").Value = 5

# Row 16
$ws.Range("A16").Value = "hammer"
$ws.Range("B16").Value = "blunt_weapons"
$ws.Range("C16").Value = 25
$ws.Range("D16").Value = 25
$ws.Range("E16").Value = 28
$ws.Range("F16").Value = 28
$ws.Range("G16").Value = 33
$ws.Range("H16").Value = 25
$ws.Range("I16").Value = 6
$ws.Range("J16").Value = 5

# Row 17
$ws.Range("A17").Value = "maul"
$ws.Range("B17").Value = "blunt_weapons"
$ws.Range("C17").Value = 25
$ws.Range("D17").Value = 25
$ws.Range("E17").Value = 30
$ws.Range("F17").Value = 34
$ws.Range("G17").Value = 38
$ws.Range("H17").Value = 27
$ws.Range("I17").Value = 7
$ws.Range("J17").Value = 5

# Row 18
$ws.Range("A18").Value = "pilum"
$ws.Range("B18").Value = "polearm_weapons"
$ws.Range("C18").Value = 30
$ws.Range("D18").Value = 30
$ws.Range("E18").Value = 25
$ws.Range("F18").Value = 18
$ws.Range("G18").Value = 15
$ws.Range("H18").Value = 3
$ws.Range("I18").Value = 3
$ws.Range("J18").Value = 4

# Row 19
$ws.Range("A19").Value = "spear"
$ws.Range("B19").Value = "polearm_weapons"
$ws.Range("C19").Value = 33
$ws.Range("D19").Value = 33
$ws.Range("E19").Value = 30
$ws.Range("F19").Value = 30
$ws.Range("G19").Value = 28
$ws.Range("H19").Value = 21
$ws.Range("I19").Value = 6
$ws.Range("J19").Value = 5

# Row 20
$ws.Range("A20").Value = "halberd"
$ws.Range("B20").Value = "polearm_weapons"
$ws.Range("C20").Value = 30
$ws.Range("D20").Value = 30
$ws.Range("E20").Value = 28
$ws.Range("F20").Value = 27
$ws.Range("G20").Value = 24
$ws.Range("H20").Value = 20
$ws.Range("I20").Value = 6
$ws.Range("J20").Value = 5

# Row 21
$ws.Range("A21").Value = "trident"
$ws.Range("B21").Value = "polearm_weapons"
$ws.Range("C21").Value = 29
$ws.Range("D21").Value = 29
$ws.Range("E21").Value = 28
$ws.Range("F21").Value = 26
$ws.Range("G21").Value = 29
$ws.Range("H21").Value = 13
$ws.Range("I21").Value = 6
$ws.Range("J21").Value = 5

# Restore the active-cell selection recorded in the saved workbook view
$ws.Range("G10").Select() | Out-Null
